$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.838080883026123
$ws.Range("B1").Value = 2.677904367446899
$ws.Range("C1").Value = 1.942438125610352
$ws.Range("D1").Value = 1.790980815887451
$ws.Range("E1").Value = 1.788223385810852
